$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Angptl3 -> Itgb3 -> ECs)
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.6698150192
$ws.Range("R2").Value = 6.0283351728
$ws.Range("S2").Value = 0.02276527781110837
$ws.Range("T2").Value = 0.02276527781110837

# Row 3 (FAPs -> Angptl3 -> Itgb3 -> FAPs)
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("S3").Value = 0.557805162587183
$ws.Range("T3").Value = 0.557805162587183

# Row 4 (FAPs -> Angptl3 -> Itgb3 -> MuSCs)
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 2.4496841966
$ws.Range("S4").Value = 0.08325842163361394
$ws.Range("T4").Value = 0.08325842163361392

# Row 5 (MuSCs -> Angptl3 -> Itgb3 -> ECs)
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 0.3392026019226667
$ws.Range("R5").Value = 3.052823417304
$ws.Range("S5").Value = 0.01152861797014227
$ws.Range("T5").Value = 0.01152861797014227

# Row 6 (MuSCs -> Angptl3 -> Itgb3 -> FAPs)
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.2824794265459323
$ws.Range("T6").Value = 0.2824794265459323

# Row 7 (MuSCs -> Angptl3 -> Itgb3 -> MuSCs)
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 1.240550345329667
$ws.Range("S7").Value = 0.04216309345202018
$ws.Range("T7").Value = 0.04216309345202016
